$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update existing row 2 (AiivoxjC / Bulgaria) in place ---
$ws.Range("G2").Value = 1.73
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.4
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("R2").Value = 1.62
$ws.Range("W2").Value = 6
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 7

# --- Step 2: insert new row at 3, fill with new Cyprus match (dKuY1awJ) ---
$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = 'dKuY1awJ'
$ws.Range("B3").Value = '25/11/2024'
$ws.Range("C3").Value = '14:00'
$ws.Range("D3").Value = 'CYPRUS - CYPRUS LEAGUE'
$ws.Range("E3").Value = 'AEK Larnaca'
$ws.Range("F3").Value = 'Karmiotissa'
$ws.Range("G3").Value = 1.24
$ws.Range("H3").Value = 5.7
$ws.Range("I3").Value = 11.5
$ws.Range("J3").Value = 1.65
$ws.Range("K3").Value = 2.62
$ws.Range("L3").Value = 8.5
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.17
$ws.Range("P3").Value = 4.65
$ws.Range("Q3").Value = 1.53
$ws.Range("R3").Value = 2.4
$ws.Range("S3").Value = 1.29
$ws.Range("T3").Value = 3.4
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 6.6
$ws.Range("Y3").Value = 9.75
$ws.Range("Z3").Value = 7.7
$ws.Range("AA3").Value = 11.25
$ws.Range("AB3").Value = 32
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 12.5
$ws.Range("AE3").Value = 27
$ws.Range("AF3").Value = 120
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 25
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 37
$ws.Range("AK3").Value = 450
$ws.Range("AL3").Value = 175
$ws.Range("AM3").Value = 120
$ws.Range("AN3").Value = 3.1
$ws.Range("AO3").Value = 5.2
$ws.Range("AP3").Value = 15.5
$ws.Range("AQ3").Value = 12.5
$ws.Range("AR3").Value = 37
$ws.Range("AS3").Value = 200
$ws.Range("AT3").Value = 3.4
$ws.Range("AU3").Value = 9.25
$ws.Range("AV3").Value = 80
$ws.Range("AW3").Value = 10.75
$ws.Range("AX3").Value = 60
$ws.Range("AY3").Value = 50
$ws.Range("AZ3").Value = 500
$ws.Range("BA3").Value = 450

# --- Step 3: update row 4 (formerly row 3, YHOF5rBm / Greece) shifted down ---
$ws.Range("G4").Value = 1.9
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 5
$ws.Range("Y4").Value = 9.5
$ws.Range("AE4").Value = 21
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 15
$ws.Range("AN4").Value = 3.75
$ws.Range("AO4").Value = 11
$ws.Range("AX4").Value = 26

# --- Step 4: update row 5 (formerly row 4, j3x6GxT7 / Romania) shifted down ---
$ws.Range("G5").Value = 2.38
$ws.Range("H5").Value = 2.8
$ws.Range("I5").Value = 3.2
$ws.Range("J5").Value = 3.25
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 2.25
$ws.Range("Z5").Value = 23
$ws.Range("AA5").Value = 23
$ws.Range("AH5").Value = 8
$ws.Range("AL5").Value = 29
$ws.Range("AO5").Value = 15
$ws.Range("AT5").Value = 2.25

# --- Step 5: insert new row at 6, fill with new Turkey Super Lig match (pAmNNiN5) ---
$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = 'pAmNNiN5'
$ws.Range("B6").Value = '25/11/2024'
$ws.Range("C6").Value = '14:00'
$ws.Range("D6").Value = 'TURKEY - SUPER LIG'
$ws.Range("E6").Value = 'Trabzonspor'
$ws.Range("F6").Value = 'Adana Demirspor'
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 6.5
$ws.Range("J6").Value = 1.83
$ws.Range("K6").Value = 2.63
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 19
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 5.5
$ws.Range("Q6").Value = 1.48
$ws.Range("R6").Value = 2.6
$ws.Range("S6").Value = 1.25
$ws.Range("T6").Value = 3.75
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("W6").Value = 9.5
$ws.Range("X6").Value = 8
$ws.Range("Y6").Value = 8.5
$ws.Range("Z6").Value = 10
$ws.Range("AA6").Value = 11
$ws.Range("AB6").Value = 21
$ws.Range("AC6").Value = 19
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 17
$ws.Range("AF6").Value = 41
$ws.Range("AG6").Value = 151
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 41
$ws.Range("AJ6").Value = 21
$ws.Range("AK6").Value = 67
$ws.Range("AL6").Value = 41
$ws.Range("AM6").Value = 41
$ws.Range("AN6").Value = 3.6
$ws.Range("AO6").Value = 6.5
$ws.Range("AP6").Value = 15
$ws.Range("AQ6").Value = 17
$ws.Range("AR6").Value = 34
$ws.Range("AS6").Value = 81
$ws.Range("AT6").Value = 3.75
$ws.Range("AU6").Value = 8
$ws.Range("AV6").Value = 51
$ws.Range("AW6").Value = 8.5
$ws.Range("AX6").Value = 34
$ws.Range("AY6").Value = 34
$ws.Range("AZ6").Value = 101
$ws.Range("BA6").Value = 101
$ws.Range("BB6").Value = 151
$ws.Range("BC6").Value = 351
$ws.Range("BD6").Value = 301

# --- Step 6: insert new row at 7, fill with new Turkey 1. Lig match (IVJXYAq9) ---
$ws.Rows("7:7").Insert()
$ws.Range("A7").Value = 'IVJXYAq9'
$ws.Range("B7").Value = '25/11/2024'
$ws.Range("C7").Value = '14:00'
$ws.Range("D7").Value = 'TURKEY - 1. LIG'
$ws.Range("E7").Value = 'Igdir FK'
$ws.Range("F7").Value = 'Amedspor'
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 2.75
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 4.75
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("Q7").Value = 2.35
$ws.Range("R7").Value = 1.57
$ws.Range("S7").Value = 1.53
$ws.Range("T7").Value = 2.38
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("W7").Value = 6
$ws.Range("X7").Value = 8.5
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 17
$ws.Range("AA7").Value = 19
$ws.Range("AB7").Value = 34
$ws.Range("AC7").Value = 7
$ws.Range("AD7").Value = 6
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 67
$ws.Range("AG7").Value = 351
$ws.Range("AH7").Value = 9.5
$ws.Range("AI7").Value = 19
$ws.Range("AJ7").Value = 15
$ws.Range("AK7").Value = 41
$ws.Range("AL7").Value = 41
$ws.Range("AM7").Value = 41
$ws.Range("AN7").Value = 3.75
$ws.Range("AO7").Value = 12
$ws.Range("AP7").Value = 26
$ws.Range("AQ7").Value = 41
$ws.Range("AR7").Value = 67
$ws.Range("AS7").Value = 201
$ws.Range("AT7").Value = 2.38
$ws.Range("AU7").Value = 9
$ws.Range("AV7").Value = 67
$ws.Range("AW7").Value = 5.5
$ws.Range("AX7").Value = 23
$ws.Range("AY7").Value = 34
$ws.Range("AZ7").Value = 81
$ws.Range("BA7").Value = 126
$ws.Range("BB7").Value = 301
$ws.Range("BC7").Value = 126
$ws.Range("BD7").Value = 126

# --- Step 7: row 8 (formerly row 5, zRQv9vQQ / Ukraine) is unchanged aside from the shift ---

Write-Host "Done. UsedRange rows:" $ws.UsedRange.Rows.Count
